$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "MARCOSVIN"
$ws.Range("B14").Value = "MARCOS"
$ws.Range("C14").Value = "VINICIUS"
$ws.Range("D14").Value = "DE JESUS LIMA"
$ws.Range("E14").Value = 23
$ws.Range("F14").Value = "RUPILOTT22@"

$ws.Range("A15").Value = "farjc97"
$ws.Range("B15").Value = "luiz"
$ws.Range("C15").Value = "fernando"
$ws.Range("D15").Value = "cunha"
$ws.Range("E15").Value = 28
$ws.Range("F15").Value = "abcdefg12@"
